# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: add a new day column AH (17-jul) with hourly prices
#  - "Gaz" sheet: append a new row (2025-07-15 / 33.35)
#  - "CO2" sheet: append a new row (2025-07-15 / 70.8)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add column AH (17-jul)
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, formatted like the other day headers (copy style from AG1)
$wsPrix.Range("AH1").Value = "17-jul"
$wsPrix.Range("AG1").Copy()
$wsPrix.Range("AH1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$prixValues = @{
    2  = 104.23
    3  = 98.38
    4  = 90.09999999999999
    5  = 83.77
    6  = 82.12
    7  = 89.16
    8  = 87.05
    9  = 101.62
    10 = 105.92
    11 = 93.09999999999999
    12 = 85
    13 = 77.06999999999999
    14 = 70.17
    15 = 56.14
    16 = 51.21
    17 = 66.56
    18 = 76.81
    19 = 84.36
    20 = 91.52
    21 = 108.66
    22 = 121.41
    23 = 125.8
    24 = 125.07
    25 = 114.18
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Cells.Item($row, 34).Value = $prixValues[$row]  # column 34 = AH
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 31 (2025-07-15 / 33.35)
# ---------------------------------------------------------------------------
# Dates in this sheet are stored as plain text (e.g. "2025-06-16"), not as
# real Excel date serials. Force the cell to text first so Excel's
# automatic date recognition doesn't turn the literal into a date serial,
# then drop the now-unneeded number format so the cell keeps the default
# (unstyled) look used by all the other rows.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A31").NumberFormat = "@"
$wsGaz.Range("A31").Value = "2025-07-15"
$wsGaz.Range("A31").ClearFormats()
$wsGaz.Range("B31").Value = 33.35

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 31 (2025-07-15 / 70.8)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A31").NumberFormat = "@"
$wsCo2.Range("A31").Value = "2025-07-15"
$wsCo2.Range("A31").ClearFormats()
$wsCo2.Range("B31").Value = 70.8
